$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RequestProcessingType" column (column F); everything to the
# right (ResponseFile, ResponseProcessingType, HTTPAction, ExcludeField,
# HttpStatusCode, StoreResponseVariables, AddifyVariables, ...) shifts left.
$ws.Range("F1").EntireColumn.Delete()

# Rename the header cells (now in their shifted positions) to the new,
# standardized field names.
$ws.Range("H1").Value = "Action"
$ws.Range("I1").Value = "ExcludeFields"
$ws.Range("J1").Value = "StatusCode"

# Restore the user's selection to K1, matching the saved workbook view.
$ws.Range("K1").Select() | Out-Null
